# 00_create_template.xlsx -- "Add files via upload"
#
# The "Data" sheet's scenario table (rows 12-21) is reshuffled:
#   - The two "Adjustment" scenario rows (formerly rows 16-17) move up to
#     become the new rows 12-13.
#   - The four "Transfer" scenario rows (formerly rows 12-15) move down to
#     become the new rows 14-17, and their G-column label is renamed from
#     "MPA_TRANSFER_SCENARIO_*" to "MPA_ADJUSTMENT_TRANSF_SCENARIO_*".
#   - The four "Retirement" rows (18-21) stay in place, but their G-column
#     label is renamed from "MPA_RETIREMENT_*" to "MPA_ADJUSTMENT_RETIRE_*"
#     (column H keeps the original "MPA_RETIREMENT_*" text).
#
# Only columns A, G, H, I and BB actually change per row; every other
# column (B, C, E, J, P, R, S, T, AQ..BD, ...) is identical across the
# affected rows and is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

function Set-Cell {
    param($addr, $value)
    $ws.Range($addr).Value2 = $value
}

# ---------------------------------------------------------------------
# 1. Write the new text/number content for the affected cells.
# ---------------------------------------------------------------------

# New rows 12-13: the former "Adjustment" rows (16-17), unchanged content,
# just relocated.
Set-Cell "A12" "11"
Set-Cell "G12" "MPA_ADJUSTMENT_SCENARIO_1"
Set-Cell "H12" "MPA_ADJUSTMENT_SCENARIO_1"
Set-Cell "I12" "SN_ADJUSTMENT1"
Set-Cell "BB12" "MANU"

Set-Cell "A13" "12"
Set-Cell "G13" "MPA_ADJUSTMENT_SCENARIO_2"
Set-Cell "H13" "MPA_ADJUSTMENT_SCENARIO_2"
Set-Cell "I13" "SN_ADJUSTMENT2"
Set-Cell "BB13" "MANU"

# New rows 14-17: the former "Transfer" rows (12-15), relocated, with the
# G-column label renamed to the "MPA_ADJUSTMENT_TRANSF_SCENARIO_*" form.
Set-Cell "A14" "7"
Set-Cell "G14" "MPA_ADJUSTMENT_TRANSF_SCENARIO_SENDER_MAIN"
Set-Cell "H14" "MPA_TRANSFER_SCENARIO_1"
Set-Cell "I14" "SN_MPATRANSFER1"
Set-Cell "BB14" "SUL3"

Set-Cell "A15" "8"
Set-Cell "G15" "MPA_ADJUSTMENT_TRANSF_SCENARIO_SENDER_SUB"
Set-Cell "H15" "MPA_TRANSFER_SCENARIO_2"
Set-Cell "I15" "SN_MPATRANSFER1"
Set-Cell "BB15" "SUL3"

Set-Cell "A16" "9"
Set-Cell "G16" "MPA_ADJUSTMENT_TRANSF_SCENARIO_INTRA_REC_MAIN"
Set-Cell "H16" "MPA_TRANSFER_SCENARIO_3"
Set-Cell "I16" "SN_MPATRANSFER1"
Set-Cell "BB16" "SUL3"

Set-Cell "A17" "10"
Set-Cell "G17" "MPA_ADJUSTMENT_TRANSF_SCENARIO_INTRA_REC_SUB"
Set-Cell "H17" "MPA_TRANSFER_SCENARIO_4"
Set-Cell "I17" "SN_MPATRANSFER1"
Set-Cell "BB17" "SUL3"

# Rows 18-21 stay put; only the G-column label text changes.
Set-Cell "G18" "MPA_ADJUSTMENT_RETIRE_PARTIAL_SCENARIO1"
Set-Cell "G19" "MPA_ADJUSTMENT_RETIRE_PARTIAL_SCENARIO2"
Set-Cell "G20" "MPA_ADJUSTMENT_RETIRE_FULL_SCENARIO1"
Set-Cell "G21" "MPA_ADJUSTMENT_RETIRE_FULLL_SCENARIO2"

# ---------------------------------------------------------------------
# 2. Fix up cell formatting so it travels with the relocated content
#    (the "Text" number format on col A/I, and the distinct "top
#    aligned" style used by BB on the Adjustment rows vs the Transfer
#    rows). Use copy/paste-special-formats from an in-row anchor cell
#    that already carries the right style, so no new style entries are
#    created.
# ---------------------------------------------------------------------

$xlPasteFormats = -4122

function Copy-Format {
    param($fromAddr, $toAddr)
    $ws.Range($fromAddr).Copy() | Out-Null
    $ws.Range($toAddr).PasteSpecial($xlPasteFormats) | Out-Null
}

# Rows 12-13 (now Adjustment): G/H need the "Text" style already used by
# column A in the same row; BB needs the style already used by column AW
# in the same row.
foreach ($r in 12, 13) {
    Copy-Format "A$r" "G$r"
    Copy-Format "A$r" "H$r"
    Copy-Format "AW$r" "BB$r"
}

# Rows 14-17 (now Transfer): G needs the style already used by column AW
# in the same row; H/BB need the style already used by column B in the
# same row.
foreach ($r in 14, 15, 16, 17) {
    Copy-Format "AW$r" "G$r"
    Copy-Format "B$r" "H$r"
    Copy-Format "B$r" "BB$r"
}

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3. Cosmetic view changes that came along with the edit.
# ---------------------------------------------------------------------

# Column G widened to fit the longer "MPA_ADJUSTMENT_TRANSF_SCENARIO_*"
# labels (36 -> 43 characters of display width).
$ws.Columns("G").ColumnWidth = 42.17

# Selection moved from G10 to D14.
$ws.Range("D14").Select()
